$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a line-break paragraph, a header table, a data table, and
#    another line-break paragraph right before the
#    "To be paid upon receipt." paragraph.
# ------------------------------------------------------------------

$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$xml = @"
<w:p $w><w:r><w:br/></w:r></w:p>
<w:tbl $w>
  <w:tblPr>
    <w:tblW w:type="auto" w:w="0"/>
    <w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="1872"/>
    <w:gridCol w:w="1872"/>
    <w:gridCol w:w="1872"/>
    <w:gridCol w:w="1872"/>
    <w:gridCol w:w="1872"/>
  </w:tblGrid>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
        <w:shd w:fill="ffffff"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:jc w:val="center"/>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
            <w:b/>
            <w:sz w:val="28"/>
          </w:rPr>
          <w:t>GARDENING SERVICES FROM 07-01-2024 TO 07-31-2024</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
        <w:shd w:fill="ffffff"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:jc w:val="center"/>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
            <w:b/>
            <w:sz w:val="28"/>
          </w:rPr>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
        <w:shd w:fill="ffffff"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:jc w:val="center"/>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
            <w:b/>
            <w:sz w:val="28"/>
          </w:rPr>
          <w:t>WORK HOURS</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
        <w:shd w:fill="ffffff"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:jc w:val="center"/>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
            <w:b/>
            <w:sz w:val="28"/>
          </w:rPr>
          <w:t>RATE</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
        <w:shd w:fill="ffffff"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:jc w:val="center"/>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
            <w:b/>
            <w:sz w:val="28"/>
          </w:rPr>
          <w:t>COST</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
<w:tbl $w>
  <w:tblPr>
    <w:tblW w:type="auto" w:w="0"/>
    <w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="1872"/>
    <w:gridCol w:w="1872"/>
    <w:gridCol w:w="1872"/>
    <w:gridCol w:w="1872"/>
    <w:gridCol w:w="1872"/>
  </w:tblGrid>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p/>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p/>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p/>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p/>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p/>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Spring cleanup and general maintenance</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>(1 ppl.)</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>3.0 h</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>$ 32/hr</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>96.0 $</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Spring cleanup and general maintenance</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>(2 ppl.)</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>2.0 h</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>$ 65/hr</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>130.0 $</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Disposal of debris &amp; weeding</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>(1 ppl.)</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>1.0 h</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>$ 32/hr</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>32.0 $</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Disposal of debris &amp; weeding</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>(2 ppl.)</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>2.0 h</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>$ 65/hr</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="1872"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>130.0 $</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
<w:p $w><w:r><w:br/></w:r></w:p>
"@

$target = $d.Content
$found = $target.Find.Execute("To be paid upon receipt.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor paragraph 'To be paid upon receipt.'"
}
$target.Collapse(1)
# Splitting off a fresh (empty) paragraph right before the anchor paragraph
# keeps the anchor paragraph itself completely untouched; the new blank
# paragraph's range is then what we fill in with our XML fragment.
$target.InsertParagraphBefore()
$target.InsertXML($xml)

# ------------------------------------------------------------------
# 2. Fix the "Yu Gothic YU Semilight" font-name typo -> "Yu Gothic UI
#    Semilight" wherever it appears (the two closing paragraphs). This
#    is a font-formatting property (w:rFonts ascii/hAnsi), not body
#    text, so it needs Font.Name rather than a text Find/Replace.
# ------------------------------------------------------------------

function Set-FontByText($searchText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0)
    if ($ok) {
        $rng.Font.Name = "Yu Gothic UI Semilight"
    }
    return $ok
}

Set-FontByText("Please do not hesitate to let me know if you have questions or concerns.") | Out-Null
Set-FontByText("Thank you!") | Out-Null

Write-Output "done"
